$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4499.5
$ws.Range("J51").Value = 4499.5
$ws.Range("L51").Value = 4499.5
$ws.Range("N51").Value = -5467.5
$ws.Range("H70").Value = 11142.429
$ws.Range("I70").Value = 12749.25
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 38247.75
$ws.Range("L70").Value = 27000
$ws.Range("M70").Value = -37977.75
$ws.Range("N70").Value = -27540
$ws.Range("H73").Value = 11142.429
$ws.Range("I73").Value = 12749.25
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 38247.75
$ws.Range("L73").Value = 27000
$ws.Range("M73").Value = -37311.75
$ws.Range("N73").Value = -28872
$ws.Range("H98").Value = 998.25
$ws.Range("J98").Value = 999
$ws.Range("L98").Value = 999
$ws.Range("N98").Value = -3995
$ws.Range("H116").Value = 3762.5
$ws.Range("J116").Value = 4323.75
$ws.Range("L116").Value = 4323.75
$ws.Range("N116").Value = -11207.75
$ws.Range("H122").Value = 998.25
$ws.Range("J122").Value = 999
$ws.Range("L122").Value = 2997
$ws.Range("N122").Value = -7897
$ws.Range("H132").Value = 1759.4445
$ws.Range("I132").Value = 1680.8695
$ws.Range("K132").Value = 5042.6085
$ws.Range("M132").Value = -2512.6085
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 961.86664
$ws.Range("I2").Value = 428.33334
$ws.Range("J2").Value = 3096
$ws.Range("K2").Value = 428.33334
$ws.Range("L2").Value = 3096
$ws.Range("M2").Value = -315.33334
$ws.Range("N2").Value = -3322
$ws.Range("H45").Value = 1749.45
$ws.Range("I45").Value = 1188.1875
$ws.Range("K45").Value = 1188.1875
$ws.Range("M45").Value = -811.1875
$ws.Range("H61").Value = 11000
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 13500
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 13500
$ws.Range("M61").Value = -788
$ws.Range("N61").Value = -13924
$ws.Range("H63").Value = 5609.6
$ws.Range("I63").Value = 3970.8572
$ws.Range("K63").Value = 3970.8572
$ws.Range("M63").Value = -3284.8572
$ws.Range("H66").Value = 5609.6
$ws.Range("I66").Value = 3970.8572
$ws.Range("K66").Value = 19854.286
$ws.Range("M66").Value = -16422.286
$ws.Range("H116").Value = 961.86664
$ws.Range("I116").Value = 428.33334
$ws.Range("J116").Value = 3096
$ws.Range("K116").Value = 428.33334
$ws.Range("L116").Value = 3096
$ws.Range("M116").Value = 1865.66666
$ws.Range("N116").Value = -7684
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 1527.4375
$ws.Range("I132").Value = 1320.5
$ws.Range("K132").Value = 3961.5
$ws.Range("M132").Value = -1431.5
$ws.Range("H136").Value = 11000
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 13500
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 40500
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -45600
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 961.86664
$ws.Range("I3").Value = 428.33334
$ws.Range("J3").Value = 3096
$ws.Range("K3").Value = 428.33334
$ws.Range("L3").Value = 3096
$ws.Range("M3").Value = -314.33334
$ws.Range("N3").Value = -3324
$ws.Range("H99").Value = 2739.1428
$ws.Range("I99").Value = 2635
$ws.Range("K99").Value = 2635
$ws.Range("M99").Value = -1137
$ws.Range("H105").Value = 2698.5
$ws.Range("I105").Value = 2698.5
$ws.Range("K105").Value = 2698.5
$ws.Range("M105").Value = -951.5
$ws.Range("H107").Value = 4466.5
$ws.Range("I107").Value = 4138
$ws.Range("K107").Value = 4138
$ws.Range("M107").Value = -2218
$ws.Range("H134").Value = 3423.0667
$ws.Range("I134").Value = 3423.0667
$ws.Range("K134").Value = 10269.2001
$ws.Range("M134").Value = -7734.2001
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 14331.667
$ws.Range("I44").Value = 2000
$ws.Range("J44").Value = 15452.728
$ws.Range("K44").Value = 2000
$ws.Range("L44").Value = 15452.728
$ws.Range("M44").Value = -1558
$ws.Range("N44").Value = -16336.728
$ws.Range("H55").Value = 38749
$ws.Range("I55").Value = 37499
$ws.Range("K55").Value = 37499
$ws.Range("M55").Value = -37184
$ws.Range("H86").Value = 4770
$ws.Range("J86").Value = 4700
$ws.Range("L86").Value = 4700
$ws.Range("N86").Value = -6946
$ws.Range("H89").Value = 4770
$ws.Range("J89").Value = 4700
$ws.Range("L89").Value = 23500
$ws.Range("N89").Value = -34732
$ws.Range("H92").Value = 33312.4
$ws.Range("J92").Value = 33312.4
$ws.Range("L92").Value = 33312.4
$ws.Range("N92").Value = -38304.4
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 226.83333
$ws.Range("I6").Value = 270.2
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = 810.5999999999999
$ws.Range("L6").Value = 30
$ws.Range("M6").Value = -697.5999999999999
$ws.Range("N6").Value = -256
$ws.Range("H68").Value = 4111.2334
$ws.Range("J68").Value = 4119.207
$ws.Range("L68").Value = 12357.621
$ws.Range("N68").Value = -13979.621
$ws.Range("H71").Value = 4111.2334
$ws.Range("J71").Value = 4119.207
$ws.Range("L71").Value = 37072.863
$ws.Range("N71").Value = -45184.863
$ws.Range("H107").Value = 1428.9
$ws.Range("J107").Value = 1485.75
$ws.Range("L107").Value = 4457.25
$ws.Range("N107").Value = -8297.25
$ws.Range("H113").Value = 1221.25
$ws.Range("I113").Value = 747
$ws.Range("J113").Value = 1379.3334
$ws.Range("K113").Value = 2241
$ws.Range("L113").Value = 4138.0002
$ws.Range("M113").Value = -71
$ws.Range("N113").Value = -8478.0002
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -15470
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 4020
$ws.Range("I133").Value = 4020
$ws.Range("K133").Value = 12060
$ws.Range("M133").Value = -7000
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 225000
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H113").Value = 1299.5
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2034
$ws.Range("I122").Value = 1975.2222
$ws.Range("K122").Value = 5925.6666
$ws.Range("M122").Value = -3475.6666
$ws.Range("H132").Value = 8473.666999999999
$ws.Range("I132").Value = 9988.5
$ws.Range("J132").Value = 5444
$ws.Range("K132").Value = 29965.5
$ws.Range("L132").Value = 16332
$ws.Range("M132").Value = -27435.5
$ws.Range("N132").Value = -21392
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 894.75
$ws.Range("J22").Value = 986.3333
$ws.Range("L22").Value = 986.3333
$ws.Range("N22").Value = -1576.3333
$ws.Range("H27").Value = 894.75
$ws.Range("J27").Value = 986.3333
$ws.Range("L27").Value = 986.3333
$ws.Range("N27").Value = -1200.3333
$ws.Range("H68").Value = 2130
$ws.Range("I68").Value = 2130
$ws.Range("K68").Value = 2130
$ws.Range("M68").Value = -1381
$ws.Range("H71").Value = 2130
$ws.Range("I71").Value = 2130
$ws.Range("K71").Value = 10650
$ws.Range("M71").Value = -6906
$ws.Range("H93").Value = 1196.8846
$ws.Range("I93").Value = 1092.3914
$ws.Range("K93").Value = 1092.3914
$ws.Range("M93").Value = 155.6086
$ws.Range("H119").Value = 65900
$ws.Range("J119").Value = 65900
$ws.Range("L119").Value = 65900
$ws.Range("N119").Value = -75576
$ws.Range("H132").Value = 3063
$ws.Range("I132").Value = 2818.8572
$ws.Range("K132").Value = 8456.571599999999
$ws.Range("M132").Value = -5926.571599999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 359060
$ws.Range("I2").Value = 15000
$ws.Range("J2").Value = 875150
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 875150
$ws.Range("M2").Value = -14888
$ws.Range("N2").Value = -875374
$ws.Range("H4").Value = 25003500
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H107").Value = 1618.2941
$ws.Range("I107").Value = 1478.1666
$ws.Range("J107").Value = 1954.6
$ws.Range("K107").Value = 4434.4998
$ws.Range("L107").Value = 5863.799999999999
$ws.Range("M107").Value = -2514.4998
$ws.Range("N107").Value = -9703.799999999999
$ws.Range("H136").Value = 2284.7
$ws.Range("I136").Value = 2284.7
$ws.Range("K136").Value = 6854.099999999999
$ws.Range("M136").Value = -4304.099999999999
